# Generate Report for Archive
# - Updates the localization status from "Ready for handoff" to "In Translation"
#   on the Overview sheet (columns E/F, row 2) and on each per-language sheet's
#   "Status" column (column C, row 2).
# - Narrows the now-shorter "Status" columns to fit the new, shorter text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Update status text: "Ready for handoff" -> "In Translation" ---
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# --- Shrink the Status columns to match the shorter text ---
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
